$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 18
$ws.Range("B2").Value = 18
$ws.Range("C2").Value = 0

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
